$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" column (D) values ---
# Force these cells to remain plain text (many values look numeric,
# e.g. "1.002", and would otherwise be auto-converted to numbers by Excel).
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.945.42"
$ws.Range("D3").Value = "1.650.47"
$ws.Range("D5").Value = "309.47"
$ws.Range("D7").Value = "0.3894"
$ws.Range("D8").Value = "0.3830"
$ws.Range("D9").Value = "51.21"
$ws.Range("D10").Value = "1.350"
$ws.Range("D11").Value = "1.002"
$ws.Range("D12").Value = "0.08435"
$ws.Range("D13").Value = "23.89"
$ws.Range("D14").Value = "7.087"
$ws.Range("D17").Value = "1.653.42"
$ws.Range("D18").Value = "94.45"
$ws.Range("D19").Value = "0.06994"
$ws.Range("D20").Value = "19.70"
$ws.Range("D21").Value = "6.933"
$ws.Range("D23").Value = "13.72"
$ws.Range("D24").Value = "23.953.64"
$ws.Range("D25").Value = "2.458"
$ws.Range("D26").Value = "2.975"
$ws.Range("D27").Value = "22.07"
$ws.Range("D28").Value = "150.94"
$ws.Range("D29").Value = "5.416"
$ws.Range("D30").Value = "138.70"
$ws.Range("D31").Value = "7.821"
$ws.Range("D32").Value = "2.492"
$ws.Range("D33").Value = "1.835.07"
$ws.Range("D34").Value = "1.048"
$ws.Range("D35").Value = "0.08077"
$ws.Range("D36").Value = "0.02962"
$ws.Range("D37").Value = "6.740"
$ws.Range("D40").Value = "0.09135"
$ws.Range("D41").Value = "0.7554"
$ws.Range("D42").Value = "13.46"
$ws.Range("D43").Value = "1.422"
$ws.Range("D44").Value = "16.43"
$ws.Range("D45").Value = "0.6948"
$ws.Range("D47").Value = "4.088"
$ws.Range("D49").Value = "0.08278"
$ws.Range("D50").Value = "134.39"
$ws.Range("D51").Value = "1.206"

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Update "Volume(1h)" column (E) values ---
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  +5.27%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +6.35%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("E38").Value = "  +5.20%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  +0.69%  "
